$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.499087929725647
$ws.Range("B1").Value = 1.648744106292725
$ws.Range("C1").Value = 1.961511969566345
$ws.Range("D1").Value = 3.017608165740967
$ws.Range("E1").Value = 1.555091500282288
